$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Don Vi" (unit) values for two rows:
#  - C8: "Ho Chi Minh" -> "HCM"
#  - C10: "De Che" -> "Ha Noi"
$ws.Range("C8").Value = "HCM"
$ws.Range("C10").Value = "Hà Nội"

# The roster only needs entries through row 10 now; drop the extra
# sample rows (11-20) and shift everything below them up, trimming the
# sheet down from 30 rows to 20.
$ws.Rows("11:20").Delete() | Out-Null

# Leave the selection on the cell that was just edited.
$ws.Range("C8").Select() | Out-Null
